$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-21 (columns A-H)
$data = @(
    @("1328206", "https://aiesec.org/opportunity/global-talent/1328206", "Power BI Specialist", "Frankfurt am Main, Deutschland", "No", "0 applicants", "3 - 6 Months", "Greyfood GmbH"),
    @("1328204", "https://aiesec.org/opportunity/global-talent/1328204", "Sales Support", "Frankfurt am Main, Deutschland", "No", "6 applicants", "3 - 6 Months", "Greyfood GmbH"),
    @("1328203", "https://aiesec.org/opportunity/global-talent/1328203", "[Impact Porto Alegre] - Business Analyst", "Santa Maria, RS, Brazil", "No", "1 applicant", "3 - 6 Months", "Sislimpa"),
    @("1328196", "https://aiesec.org/opportunity/global-talent/1328196", "Marketing Trainee Specialties", "Panamá, Provincia de Panamá, Panamá", "No", "3 applicants", "6 - 18 Months", "Merck"),
    @("1328194", "https://aiesec.org/opportunity/global-talent/1328194", "Junior Buying Specialist", "Porto, Portugal", "No", "9 applicants", "3 - 6 Months", "SDSR - Sports Division SR"),
    @("1328185", "https://aiesec.org/opportunity/global-talent/1328185", "Data Scientist", "Frankfurt am Main, Deutschland", "No", "18 applicants", "3 - 6 Months", "Greyfood GmbH"),
    @("1328156", "https://aiesec.org/opportunity/global-talent/1328156", "Young Professional Comercial", "San Salvador, El Salvador", "No", "5 applicants", "6 - 18 Months", "Samsung Electronics Latinoamerica"),
    @("1328131", "https://aiesec.org/opportunity/global-talent/1328131", "Content Creation Intern", "Panchkula, Haryana, India", "No", "1 applicant", "3 - 6 Months", "Idea Clan"),
    @("1328130", "https://aiesec.org/opportunity/global-talent/1328130", "Digital Marketing Intern", "Sahibzada Ajit Singh Nagar, Punjab, India", "No", "2 applicants", "3 - 6 Months", "Requisite Technologies Pvt Ltd"),
    @("1328115", "https://aiesec.org/opportunity/global-talent/1328115", "Business Development Executive", "Sahibzada Ajit Singh Nagar, Punjab, India", "No", "0 applicants", "3 - 6 Months", "Requisite Technologies Pvt Ltd"),
    @("1328051", "https://aiesec.org/opportunity/global-talent/1328051", "AI Intern", "Sahibzada Ajit Singh Nagar, Punjab, India", "No", "3 applicants", "3 - 6 Months", "Requisite Technologies Pvt Ltd"),
    @("1327881", "https://aiesec.org/opportunity/global-talent/1327881", "Preschool Coordinator", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "The Hungry Mind Academy"),
    @("1326913", "https://aiesec.org/opportunity/global-talent/1326913", "Occupational Health and Safety Project Specialist (Mine Opp)", "Mexico City, CDMX, Mexico", "No", "21 applicants", "6 - 18 Months", "Sodexo Mexico"),
    @("1326890", "https://aiesec.org/opportunity/global-talent/1326890", "Occupational Health and Safety Projects Specialist", "Mexico City, CDMX, Mexico", "No", "14 applicants", "6 - 18 Months", "Sodexo Mexico"),
    @("1326590", "https://aiesec.org/opportunity/global-talent/1326590", "BUSINESS DEVELOPMENT", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "28 applicants", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326575", "https://aiesec.org/opportunity/global-talent/1326575", "Graphic design", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "9 applicants", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326555", "https://aiesec.org/opportunity/global-talent/1326555", "Graphic desgin", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "1 applicant", "9 - 12 Weeks", "Egypt holiday travel"),
    @("1326467", "https://aiesec.org/opportunity/global-talent/1326467", "Customer Relations Executive/ Intern", "Colombo, Sri Lanka", "No", "72 applicants", "6 - 18 Months", "Formix (Pvt) Ltd"),
    @("1324500", "https://aiesec.org/opportunity/global-talent/1324500", "DT Software Engineer Trainee (EU ONLY)", "Brussels, Belgium", "No", "63 applicants", "6 - 18 Months", "UCB"),
    @("1321275", "https://aiesec.org/opportunity/global-talent/1321275", "Business Development Executive", "Cairo, Cairo Governorate, Egypt", "No", "22 applicants", "3 - 6 Months", "Silverkey Technologies Egypt"),
)

# Column A holds numeric-looking opportunity IDs that must stay text,
# exactly like the source data (strings, not numbers).
$ws.Range("A2:A21").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Drop the temporary text number-format now that the IDs are locked in as text.
$ws.Range("A2:A21").Style = "Normal"

# Clear the old premium (yellow fill) highlighting in the PREMIUM column --
# no opportunity in the refreshed scrape is flagged premium any more.
$ws.Range("E2:E21").Style = "Normal"

# Column width adjustments (C: 64->63, D: 34->70, H: 33->36)
$ws.Columns.Item(3).ColumnWidth = 62.166666666666664
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667
$ws.Columns.Item(8).ColumnWidth = 35.166666666666664
